$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: add new "Ref" value (B13) and update the MVC description (E13) ---
$ws.Range("B13").Value = "01"
$ws.Range("E13").Value = "Model View Controller (MVC)"

# --- Insert three fresh rows (14:16) for the new Use case / architecture scenarios ---
$ws.Range("A14:A16").EntireRow.Insert()

# Row 14: A14 gets the check-mark style used by A5:A9 (copy formats from A5), plus value
$ws.Range("A5").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = [char]0x221A
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("E14").Clear()
$ws.Range("D14").Value = "  -Application Architecture"

# Row 15
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()
$ws.Range("E15").Clear()
$ws.Range("D15").Value = "  -Server Architecture"

# Row 16
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()
$ws.Range("E16").Clear()
$ws.Range("D16").Value = "  -Client Architecture"

# --- Update the selected cell to match the author's final cursor position ---
$ws.Range("E12").Select()

Write-Output "done"
